# YumaPointImport.xlsx edit
#
# Author intent (from commit message "Correct Tests ... Point work."):
# the "Longitude" calculated column in Table1 (cell AB2, on the single
# data row) was wrongly copied from the Latitude column's formula
# (=J2) - it needs to point at the Longitude source column instead
# (=K2). Fix the formula so the table recalculates the correct
# longitude value (-112.229061 instead of the latitude 36.0799823).
#
# The workbook's last-used selection/scroll position is also updated
# to reflect where the author was working (around the newly-fixed
# Longitude column) when the file was saved.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Core fix: Longitude column formula referenced the wrong source
# column (Latitude's J2) -- point it at Longitude's own source, K2.
$ws.Range("AB2").Formula = "=K2"

# --- Window/selection bookkeeping: the sheet was left scrolled to
# column W with the active cell on AB3 (just under the corrected
# Longitude cell).
$ws.Range("AB3").Select()
$win = $excel.ActiveWindow
$win.ScrollColumn = 23
$win.ScrollRow = 1
